$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet AVIO.MI -> TGYM.MI
$ws.Name = "TGYM.MI"

# 2. Insert a new column before T ("issuance_of_debt"), shifting T..AF to U..AG
$ws.Range("T1").EntireColumn.Insert()
$ws.Range("T1").Value = "issuance_of_debt"

# 3. Update data rows 2-5 with the new TGYM.MI figures
$ws.Range("A2").Value = "TGYM.MI"
$ws.Range("B2").Value = "balance_sheet"
$ws.Range("C2").Value = 44561
$ws.Range("D2").Value = 309841000
$ws.Range("E2").Value = 388938000
$ws.Range("F2").Value = 259145000
$ws.Range("G2").Value = 201327500
$ws.Range("H2").Value = 763092000
$ws.Range("I2").Value = 322341000
$ws.Range("J2").Value = 103130000
$ws.Range("K2").Value = 451531000
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = 124477000
$ws.Range("N2").Value = "cash_flow"
$ws.Range("O2").Value = -27691000
$ws.Range("P2").Value = 174306000
$ws.Range("Q2").Value = -127266000
$ws.Range("R2").Value = 65516000
$ws.Range("S2").Value = 4323000
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 93207000
$ws.Range("V2").Value = -25000000
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = "income_stmt"
$ws.Range("Y2").Value = 0.31
$ws.Range("Z2").Value = 80598000
$ws.Range("AA2").Value = 117004000
$ws.Range("AB2").Value = 243759000
$ws.Range("AC2").Value = 63065000
$ws.Range("AD2").Value = 64185000
$ws.Range("AE2").Value = 79731000
$ws.Range("AF2").Value = 16466000
$ws.Range("AG2").Value = 611412000
$ws.Range("A3").Value = "TGYM.MI"
$ws.Range("B3").Value = "balance_sheet"
$ws.Range("C3").Value = 44926
$ws.Range("D3").Value = 343922000
$ws.Range("E3").Value = 375367000
$ws.Range("F3").Value = 288439000
$ws.Range("G3").Value = 201327500
$ws.Range("H3").Value = 786207000
$ws.Range("I3").Value = 346547000
$ws.Range("J3").Value = 69220000
$ws.Range("K3").Value = 440280000
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = 141471000
$ws.Range("N3").Value = "cash_flow"
$ws.Range("O3").Value = -34963000
$ws.Range("P3").Value = 205358000
$ws.Range("Q3").Value = -50795000
$ws.Range("R3").Value = 77280000
$ws.Range("S3").Value = -32112000
$ws.Range("T3").Value = 4006000
$ws.Range("U3").Value = 112243000
$ws.Range("V3").Value = -54848000
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = "income_stmt"
$ws.Range("Y3").Value = 0.32
$ws.Range("Z3").Value = 84140000
$ws.Range("AA3").Value = 124504000
$ws.Range("AB3").Value = 278106000
$ws.Range("AC3").Value = 63587000
$ws.Range("AD3").Value = 82631000
$ws.Range("AE3").Value = 83309000
$ws.Range("AF3").Value = 19434000
$ws.Range("AG3").Value = 721490000
$ws.Range("A4").Value = "TGYM.MI"
$ws.Range("B4").Value = "balance_sheet"
$ws.Range("C4").Value = 45291
$ws.Range("D4").Value = 354544000
$ws.Range("E4").Value = 376629000
$ws.Range("F4").Value = 297859000
$ws.Range("G4").Value = 200490528
$ws.Range("H4").Value = 818309000
$ws.Range("I4").Value = 354544000
$ws.Range("J4").Value = 65900000
$ws.Range("K4").Value = 454597000
$ws.Range("L4").Value = 836972
$ws.Range("M4").Value = 158171000
$ws.Range("N4").Value = "cash_flow"
$ws.Range("O4").Value = -36179000
$ws.Range("P4").Value = 224730000
$ws.Range("Q4").Value = -50936000
$ws.Range("R4").Value = 70642000
$ws.Range("S4").Value = -30819000
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = 106821000
$ws.Range("V4").Value = -11735000
$ws.Range("W4").Value = -6922000
$ws.Range("X4").Value = "income_stmt"
$ws.Range("Y4").Value = 0.37
$ws.Range("Z4").Value = 101142000
$ws.Range("AA4").Value = 146709000
$ws.Range("AB4").Value = 324677000
$ws.Range("AC4").Value = 73640000
$ws.Range("AD4").Value = 94699000
$ws.Range("AE4").Value = 100400000
$ws.Range("AF4").Value = 23232000
$ws.Range("AG4").Value = 808091000
$ws.Range("A5").Value = "TGYM.MI"
$ws.Range("B5").Value = "balance_sheet"
$ws.Range("C5").Value = 45657
$ws.Range("D5").Value = 378996000
$ws.Range("E5").Value = 402722000
$ws.Range("F5").Value = 326301000
$ws.Range("G5").Value = 199161715
$ws.Range("H5").Value = 904134000
$ws.Range("I5").Value = 378996000
$ws.Range("J5").Value = 70817000
$ws.Range("K5").Value = 517340000
$ws.Range("L5").Value = 2165785
$ws.Range("M5").Value = 173773000
$ws.Range("N5").Value = "cash_flow"
$ws.Range("O5").Value = -41556000
$ws.Range("P5").Value = 268709000
$ws.Range("Q5").Value = -73019000
$ws.Range("R5").Value = 113852000
$ws.Range("S5").Value = -38611000
$ws.Range("T5").Value = ""
$ws.Range("U5").Value = 155408000
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = -13128000
$ws.Range("X5").Value = "income_stmt"
$ws.Range("Y5").Value = ""
$ws.Range("Z5").Value = 123883000
$ws.Range("AA5").Value = 175183000
$ws.Range("AB5").Value = 373496000
$ws.Range("AC5").Value = 87041000
$ws.Range("AD5").Value = 120940000
$ws.Range("AE5").Value = 122962000
$ws.Range("AF5").Value = 33846000
$ws.Range("AG5").Value = 901289000
